# Insert a new row for a new item ("مبرد اظافر") above the existing
# "مناديل جيب مبلله" row, pushing it (and the totals/footer rows) down by one,
# and bump the footer timestamp by one minute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 currently holds item #14 "مناديل جيب مبلله". Insert a new blank row
# above it; everything from row 20 downward shifts down by one.
$ws.Rows(20).Insert()

# Copy the (now shifted-down) item row's formatting into the freshly inserted
# row so the new row matches the sheet's existing item-row style.
$ws.Range("A21:Q21").Copy()
$ws.Range("A20:Q20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A20").PasteSpecial(-4122)

# New item data in row 20 ("مبرد اظافر" / nail file).
$ws.Range("A20").Value = 14
$ws.Range("C20").Value = "مبرد اظافر"
$ws.Range("H20").Value = "1:0"
$ws.Range("L20").Value = "0"
$ws.Range("N20").Value = "15.00"
$ws.Range("P20").Value = "15.0000"
$ws.Range("Q20").Value = "1:0"

# Renumber the pushed-down item row (was #14, now #15) - its other data is
# unchanged.
$ws.Range("A21").Value = 15

# Update the running total (was 540.67, now +15.00 for the new item).
$ws.Range("P22").Value = 555.67

# Update the footer timestamp by one minute.
$ws.Range("A23").Value = "Sunday, 17 August, 2025 10:03 AM"

Write-Output "done"
